# 4_DefinitionenDerProzesswörter.xlsx - add two new process words
# "einlesen" and "zulassen" as new rows at the bottom of the list
# (Tabelle1, column A), matching the existing formatting used by the
# other single-column entries (e.g. "pausieren", "fortsetzen", ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Append the two new process words in column A, rows 26 and 27.
$ws.Range("A26").Value = "einlesen"
$ws.Range("A27").Value = "zulassen"

# Match formatting of the preceding single-column rows (wrap text,
# top-aligned), which corresponds to the shared cell style already
# used throughout column A.
$ws.Range("A26:A27").WrapText = $true
$ws.Range("A26:A27").VerticalAlignment = -4160

# Reflect the author's last on-screen scroll position / selection.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B28").Select() | Out-Null
